$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Login - invalid credentials): Stop Iteration count changes from 5 to 4
# Leading apostrophe preserves the "quote-prefixed text" storage (numeric-looking text)
# so the cell keeps its original text-style rather than becoming a true number.
$ws.Range("F3").Value = "'4"

# Row 6 (Verify SignUp with valid credentials): now executed -> Execute flag flips to "Yes"
$ws.Range("B6").Value = "Yes"

# Row 6 H6: SignUp test-data keyword string gains the new verification / PIN / Face ID params
$ws.Range("H6").Value = "coyni_mobile.tests.SignUpTest,
testSignUp,
-pcreateAccount,
-pfirstName,
-plastName,
-pemail,
-pphoneNumber,
-ppassword,
-pconfirmPassword,
-pphoneVerificationHeading,
-pcode,
-pemailVerificationHeading,
-psecureYourAccountHeading,
-pchoosePinHeading,
-ppin,
-pconfirmPinHeading,
-penableFaceIdHeading,
-pcreateAccountHeading"

# Row 8 (Verify field validation in create account): no longer executed -> Execute flag flips to "No"
$ws.Range("B8").Value = "No"

# Reflect the author's final on-screen selection (cell H6 was the one edited last).
$ws.Range("H6").Select() | Out-Null
